$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2").Value = 0.9717261791229248
$ws.Range("F3").Value = 0.9714771509170532
$ws.Range("F4").Value = 0.9699634909629822
$ws.Range("F5").Value = 0.6271273493766785
$ws.Range("F6").Value = 0.973065972328186
$ws.Range("F7").Value = 0.9740249514579773
$ws.Range("F8").Value = 0.9738725423812866
$ws.Range("F9").Value = 0.9732975363731384
$ws.Range("F10").Value = 0.9738725423812866
$ws.Range("F11").Value = 0.5309058427810669
$ws.Range("F12").Value = 0.975741982460022
$ws.Range("F13").Value = 0.8982094526290894
$ws.Range("F14").Value = 0.9727062582969666
$ws.Range("F15").Value = 0.9727745056152344
$ws.Range("F16").Value = 0.9730832576751709
$ws.Range("F17").Value = 0.9731369614601135
$ws.Range("F18").Value = 0.2995750606060028
$ws.Range("F19").Value = 0.9745760560035706
$ws.Range("F20").Value = 0.9749926924705505
$ws.Range("F21").Value = 0.9702966809272766
$ws.Range("F22").Value = 0.9728440046310425
$ws.Range("F23").Value = 0.9745829701423645
$ws.Range("F24").Value = 0.9750985503196716
$ws.Range("F25").Value = 0.6264590620994568
$ws.Range("F26").Value = 0.291460245847702
$ws.Range("F27").Value = 0.9717508554458618
$ws.Range("F28").Value = 0.9714771509170532
$ws.Range("F29").Value = 0.97041255235672
$ws.Range("F30").Value = 0.9730724096298218
$ws.Range("F31").Value = 0.9715277552604675
$ws.Range("F32").Value = 0.9741558432579041
$ws.Range("F33").Value = 0.8576926589012146
$ws.Range("F34").Value = 0.8988063335418701
$ws.Range("F35").Value = 0.8929467797279358
$ws.Range("F36").Value = 0.974073588848114
$ws.Range("F37").Value = 0.8947763442993164
$ws.Range("F38").Value = 0.9754317998886108
$ws.Range("F39").Value = 0.8945091962814331
$ws.Range("F40").Value = 0.8572773337364197
$ws.Range("F41").Value = 0.975463330745697
$ws.Range("F42").Value = 0.8593934178352356
$ws.Range("F43").Value = 0.8626639246940613
$ws.Range("F44").Value = 0.8958331942558289
$ws.Range("F45").Value = 0.8655011653900146
$ws.Range("F46").Value = 0.8561663031578064
$ws.Range("F47").Value = 0.6203930974006653
$ws.Range("F48").Value = 0.8539775013923645
$ws.Range("F49").Value = 0.9751623868942261
$ws.Range("F50").Value = 0.8592694997787476
$ws.Range("F51").Value = 0.8541398048400879
$ws.Range("F52").Value = 0.9712916612625122
$ws.Range("F53").Value = 0.8592694997787476
$ws.Range("F54").Value = 0.8595391511917114
$ws.Range("F55").Value = 0.9739654064178467
$ws.Range("F56").Value = 0.9742932319641113
$ws.Range("F57").Value = 0.9756810665130615
$ws.Range("F58").Value = 0.892865777015686
$ws.Range("F59").Value = 0.975336492061615
$ws.Range("F60").Value = 0.9734772443771362
$ws.Range("F61").Value = 0.9741453528404236
$ws.Range("F62").Value = 0.9752413630485535
$ws.Range("F63").Value = 0.8959683179855347
$ws.Range("F64").Value = 0.469641238451004
$ws.Range("F65").Value = 0.8929812908172607
$ws.Range("F66").Value = 0.8929812908172607
$ws.Range("F67").Value = 0.8959721922874451
$ws.Range("F68").Value = 0.8929812908172607
$ws.Range("F69").Value = 0.644463062286377
$ws.Range("F70").Value = 0.9740021228790283
$ws.Range("F71").Value = 0.895367443561554
$ws.Range("F72").Value = 0.9685226678848267
$ws.Range("F73").Value = 0.8947239518165588
$ws.Range("F74").Value = 0.972861111164093
$ws.Range("F75").Value = 0.973066508769989
$ws.Range("F76").Value = 0.8947239518165588
$ws.Range("F77").Value = 0.9715393781661987
$ws.Range("F78").Value = 0.8947239518165588
$ws.Range("F79").Value = 0.9745737910270691
$ws.Range("F80").Value = 0.8947955369949341
$ws.Range("F81").Value = 0.971558690071106
$ws.Range("F82").Value = 0.8947955369949341
$ws.Range("F83").Value = 0.9644967317581177
$ws.Range("F84").Value = 0.9644967317581177
$ws.Range("F85").Value = 0.9756811261177063
$ws.Range("F86").Value = 0.8957412838935852
$ws.Range("F87").Value = 0.9745206832885742
$ws.Range("F88").Value = 0.8947451710700989
$ws.Range("F89").Value = 0.9715264439582825
$ws.Range("F90").Value = 0.972831666469574
$ws.Range("F91").Value = 0.9728267192840576
$ws.Range("F92").Value = 0.97284996509552
$ws.Range("F93").Value = 0.9728230237960815
$ws.Range("F94").Value = 0.9728220105171204
$ws.Range("F95").Value = 0.9728443026542664
$ws.Range("F96").Value = 0.9758037328720093
$ws.Range("F97").Value = 0.2213776111602783
$ws.Range("F98").Value = 0.9750561714172363
$ws.Range("F99").Value = 0.971674382686615
$ws.Range("F100").Value = 0.9703823328018188
$ws.Range("F101").Value = 0.9746648669242859
$ws.Range("F102").Value = 0.6244308948516846
$ws.Range("F103").Value = 0.8594067096710205
$ws.Range("F104").Value = 0.9730537533760071
$ws.Range("F105").Value = 0.6553149223327637
$ws.Range("F106").Value = 0.9704095125198364
$ws.Range("F107").Value = 0.8636355400085449
$ws.Range("F108").Value = 0.6099866032600403
$ws.Range("F109").Value = 0.9717712998390198
$ws.Range("F110").Value = 0.9714771509170532
$ws.Range("F111").Value = 0.8577429056167603
$ws.Range("F112").Value = 0.8626078963279724
$ws.Range("F113").Value = 0.8971216678619385
$ws.Range("F114").Value = 0.9752607941627502
$ws.Range("F115").Value = 0.8970744013786316
$ws.Range("F116").Value = 0.96440589427948
$ws.Range("F117").Value = 0.9681317806243896
$ws.Range("F118").Value = 0.6112022399902344
$ws.Range("F119").Value = 0.9738143086433411
$ws.Range("F120").Value = 0.9714771509170532
$ws.Range("F121").Value = 0.6452741622924805
$ws.Range("F122").Value = 0.9751656651496887
$ws.Range("F123").Value = 0.6452741622924805
$ws.Range("F124").Value = 0.9751656651496887
$ws.Range("F125").Value = 0.6452741622924805
$ws.Range("F126").Value = 0.9748996496200562
$ws.Range("F127").Value = 0.9717430472373962
$ws.Range("F128").Value = 0.9714771509170532
$ws.Range("F129").Value = 0.9738816022872925
$ws.Range("F130").Value = 0.9702634811401367
$ws.Range("F131").Value = 0.9680851697921753
$ws.Range("F132").Value = 0.8477449417114258
$ws.Range("F133").Value = 0.9731167554855347
$ws.Range("F134").Value = 0.9727917313575745
$ws.Range("F135").Value = 0.9740479588508606
$ws.Range("F136").Value = 0.9638451933860779
$ws.Range("F137").Value = 0.9730864763259888
$ws.Range("F138").Value = 0.9681291580200195
$ws.Range("F139").Value = 0.9681397080421448
$ws.Range("F140").Value = 0.9714975357055664
$ws.Range("F141").Value = 0.5270276665687561
$ws.Range("F142").Value = 0.9714771509170532
$ws.Range("F143").Value = 0.8929688930511475
$ws.Range("F144").Value = 0.9734880924224854
$ws.Range("F145").Value = 0.8929688930511475
$ws.Range("F146").Value = 0.9643860459327698
$ws.Range("F147").Value = 0.644463062286377
$ws.Range("F148").Value = 0.9716962575912476
$ws.Range("F149").Value = 0.8946552872657776
$ws.Range("F150").Value = 0.6472999453544617
$ws.Range("F151").Value = 0.9714771509170532
$ws.Range("F152").Value = 0.6002638339996338
$ws.Range("F153").Value = 0.9754392504692078
$ws.Range("F154").Value = 0.6143796443939209
$ws.Range("F155").Value = 0.9750345349311829
$ws.Range("F156").Value = 0.6147691607475281
$ws.Range("F157").Value = 0.9734753966331482
$ws.Range("F158").Value = 0.6148035526275635
$ws.Range("F159").Value = 0.9738812446594238
$ws.Range("F160").Value = 0.8963732719421387
$ws.Range("F161").Value = 0.9718206524848938
$ws.Range("F162").Value = 0.9649969935417175
$ws.Range("F163").Value = 0.9681324362754822
$ws.Range("F164").Value = 0.8928683400154114
$ws.Range("F165").Value = 0.9757151007652283
$ws.Range("F166").Value = 0.9645726680755615
$ws.Range("F167").Value = 0.8499478697776794
$ws.Range("F168").Value = 0.8968327045440674
$ws.Range("F169").Value = 0.9714771509170532
$ws.Range("F170").Value = 0.8499478697776794
$ws.Range("F171").Value = 0.8957864046096802
$ws.Range("F172").Value = 0.9714771509170532
$ws.Range("F173").Value = 0.968089759349823
$ws.Range("F174").Value = 0.9681097269058228
$ws.Range("F175").Value = 0.968089759349823
$ws.Range("F176").Value = 0.968089759349823
$ws.Range("F177").Value = 0.9681097269058228
$ws.Range("F178").Value = 0.9717199802398682
$ws.Range("F179").Value = 0.9685471057891846
$ws.Range("F180").Value = 0.9730931520462036
$ws.Range("F181").Value = 0.9748011231422424
$ws.Range("F182").Value = 0.9740332365036011
$ws.Range("F183").Value = 0.6531647443771362
$ws.Range("F184").Value = 0.9751391410827637
$ws.Range("F185").Value = 0.9754480123519897
$ws.Range("F186").Value = 0.9743496775627136
$ws.Range("F187").Value = 0.9738466739654541
$ws.Range("F188").Value = 0.9717810153961182
$ws.Range("F189").Value = 0.9727993011474609
$ws.Range("F190").Value = 0.9728150963783264
$ws.Range("F191").Value = 0.2115035057067871
$ws.Range("F192").Value = 0.9740918278694153
$ws.Range("F193").Value = 0.973141610622406
$ws.Range("F194").Value = 0.6453492045402527
$ws.Range("F195").Value = 0.8590968251228333
$ws.Range("F196").Value = 0.975081741809845
$ws.Range("F197").Value = 0.9738476872444153
$ws.Range("F198").Value = 0.9742869734764099
$ws.Range("F199").Value = 0.974597692489624
$ws.Range("F200").Value = 0.8944209814071655
$ws.Range("F201").Value = 0.974617600440979
$ws.Range("F202").Value = 0.9746358394622803
$ws.Range("F203").Value = 0.973113477230072
$ws.Range("F204").Value = 0.8501247763633728
$ws.Range("F205").Value = 0.9747381210327148
$ws.Range("F206").Value = 0.9680836200714111
$ws.Range("F207").Value = 0.2211202681064606
$ws.Range("F208").Value = 0.9685140252113342
$ws.Range("F209").Value = 0.3089878857135773
$ws.Range("F210").Value = 0.9730339050292969
